# Update "想去人数" (F column) values across all four sheets to reflect
# refreshed counts from the data source (output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1580
$ws.Range("F5").Value  = 9016
$ws.Range("F6").Value  = 256
$ws.Range("F7").Value  = 112
$ws.Range("F8").Value  = 1261
$ws.Range("F10").Value = 580
$ws.Range("F13").Value = 131
$ws.Range("F14").Value = 287
$ws.Range("F17").Value = 1475
$ws.Range("F18").Value = 1316
$ws.Range("F21").Value = 1358
$ws.Range("F22").Value = 75
$ws.Range("F25").Value = 85
$ws.Range("F26").Value = 48
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 294
$ws.Range("F30").Value = 1064
$ws.Range("F32").Value = 28
$ws.Range("F33").Value = 222
$ws.Range("F34").Value = 190
$ws.Range("F39").Value = 129
$ws.Range("F41").Value = 152
$ws.Range("F42").Value = 11
$ws.Range("F43").Value = 484
$ws.Range("F44").Value = 1226
$ws.Range("F45").Value = 682
$ws.Range("F46").Value = 205
$ws.Range("F47").Value = 44

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 222
$ws.Range("F13").Value = 158
$ws.Range("F16").Value = 668
$ws.Range("F20").Value = 70
$ws.Range("F26").Value = 219
$ws.Range("F29").Value = 201
$ws.Range("F39").Value = 95

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 744
$ws.Range("F6").Value = 292
$ws.Range("F8").Value = 2035
$ws.Range("F9").Value = 3063

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1580
$ws.Range("F5").Value  = 744
$ws.Range("F6").Value  = 9016
$ws.Range("F7").Value  = 292
$ws.Range("F10").Value = 256
$ws.Range("F11").Value = 2035
$ws.Range("F12").Value = 3063
$ws.Range("F13").Value = 222
$ws.Range("F14").Value = 1261
$ws.Range("F15").Value = 158
$ws.Range("F16").Value = 580
$ws.Range("F18").Value = 287
$ws.Range("F20").Value = 1475
$ws.Range("F21").Value = 1316
$ws.Range("F23").Value = 1358
$ws.Range("F24").Value = 75
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 48
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 1064
$ws.Range("F31").Value = 70
$ws.Range("F32").Value = 222
$ws.Range("F34").Value = 219
$ws.Range("F40").Value = 201
$ws.Range("F42").Value = 484
$ws.Range("F43").Value = 682
$ws.Range("F46").Value = 205
$ws.Range("F48").Value = 95
